$d = $word.ActiveDocument

# 1. Insert the new sentence about Cholesky decomposition / sim_distribution
#    right before "Another step ..." (this keeps the whole insertion within
#    the single italic run for now; it gets split into separate runs below,
#    mirroring how "data.table" is already split out earlier in this same
#    paragraph).
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "speed up the code. Another step",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "speed up the code. In other parts of the package, we learned that there is little we can do to reduce the processing time (e.g. Cholesky decomposition is a factor limiting the speed of the sim_distribution function). Another step",
    2)
if (-not $found1) {
    Write-Output "WARNING: could not find text for edit #1 (Cholesky/sim_distribution insertion)"
}

# 2. "many" -> "some" of the core functions
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "translate many of the core functions",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "translate some of the core functions",
    2)
if (-not $found2) {
    Write-Output "WARNING: could not find text for edit #2 (many -> some)"
}

# 3. Locate the newly-inserted "sim_distribution" and give it the
#    VerbatimChar character style (as used elsewhere in the letter for
#    code/package identifiers like "data.table"). Force the adjoining
#    single-space runs to split off into their own runs too, matching the
#    existing "data.table" run pattern, by toggling italic off then back on
#    (a no-op in value, but it forces a run boundary at that point).
$r3 = $d.Content
$found3 = $r3.Find.Execute("sim_distribution", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)

if ($found3) {
    $spaceBefore = $d.Range($r3.Start - 1, $r3.Start)
    $spaceBefore.Font.Italic = 0
    $spaceBefore.Font.Italic = 1

    $spaceAfter = $d.Range($r3.End, $r3.End + 1)
    $spaceAfter.Font.Italic = 0
    $spaceAfter.Font.Italic = 1

    $r3.Style = "VerbatimChar"
} else {
    Write-Output "WARNING: could not find 'sim_distribution' to style"
}
